# "Scripts created for PDP"
#
# 1) Two existing Run Flags on AppControl flip from Y to N (Cart page suite
#    and Order Summary Page suite are turned off for this run).
# 2) Three new PDP (Product Description Page) test cases are appended to the
#    Suite6 worksheet, each flagged to run ("Y") and marked "Done".
# 3) The workbook is left with the Suite6 tab active/selected, matching the
#    author's last position when the scripts were added.

$wb = $excel.ActiveWorkbook

$appControl = $wb.Worksheets.Item("AppControl")
$appControl.Range("B3").Value = "N"
$appControl.Range("B6").Value = "N"

$suite6 = $wb.Worksheets.Item("Suite6")

$suite6.Range("A7").Value = "IGP_TC_108"
$suite6.Range("B7").Value = "Y"
$suite6.Range("C7").Value = "Product name and Price-:Ensure that the selected product name and price   should be same with what we are selected in listing page."
$suite6.Range("D7").Value = "Done"

$suite6.Range("A8").Value = "IGP_TC_109"
$suite6.Range("B8").Value = "Y"
$suite6.Range("C8").Value = "Size option-: Ensure that the Size option displayed when size variants exist."
$suite6.Range("D8").Value = "Done"

$suite6.Range("A9").Value = "IGP_TC_110"
$suite6.Range("B9").Value = "Y"
$suite6.Range("C9").Value = "Need Help?-:Ensure that user able to send query by clicking on 'Need Help' from the `"Product Description`" page."
$suite6.Range("D9").Value = "Done"

# Row heights follow the wrapped "Desc" text, same as the pre-existing rows.
$suite6.Rows.Item(6).RowHeight = 28.35
$suite6.Rows.Item(7).RowHeight = 41.75
$suite6.Rows.Item(8).RowHeight = 28.35
$suite6.Rows.Item(9).RowHeight = 41.75

# Leave the cursor/selection the way the author last left it, and make
# Suite6 the active tab.
$appControl.Range("B8").Select()
$suite6.Activate()
$suite6.Range("C10").Select()
